# Fruta / hortaliza, semanal
# Insert a new weekly record for "Vega Modelo de Temuco" / Papaya at row 25,
# pushing the existing rows 25-59 down to 26-60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 25..59 down by one row (mirrors the canonical diff: a brand new
# row appears at 25 and the previous rows 25-59 become 26-60).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 44483
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100108
$ws.Range("H25").Value = "Tropicales y subtropicales"
$ws.Range("I25").Value = 100108004
$ws.Range("J25").Value = "Papaya"
$ws.Range("K25").Value = "Cultivar IV Región"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 25
$ws.Range("N25").Value = 20000
$ws.Range("O25").Value = 20000
$ws.Range("P25").Value = 20000
$ws.Range("Q25").Value = "$/bandeja 10 kilos"
$ws.Range("R25").Value = "Provincia del Elquí"
$ws.Range("S25").Value = 2000
$ws.Range("T25").Value = 10
